$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the paired "_old" / "_new" column headers (A1:J1 and L1:U1) to
# carry the actual format-version suffixes ("_FV2210" / "_FV2304")
# instead of the generic "_old" / "_new" markers. Column K ("diff")
# keeps its name - it separates the two blocks.
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($baseNames[$i] + "_FV2210")
}
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($baseNames[$i] + "_FV2304")
}

# Turn the used range into a real Excel Table (ListObject) with the
# (now renamed) header row, matching the full data range A1:U75.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U75"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split/frozen pane below row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
